# edit.ps1
# Implements the TC01 -> TC02 "single case add to cart" content edit:
#  - A2: QA site URL text/hyperlink updated (bento-qa.bento-tools.org), hyperlink re-pointed with
#        no stale in-sheet #/ location/display text (matches a plain re-insert of the link)
#  - B2: Neo4j query text updated (er_status filter Negative -> Positive, endocrine_therapy_type /
#        head(labels(samp)) filters removed)
#  - C2/D2: expected-output file names updated from the TC01 sample to the TC02 sample
#  - Row 2 height / view nudged to match the refreshed layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B2: updated Neo4j query text ---
$newQuery = @'
MATCH (ss:study_subject)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)<-[:sample_of_study_subject]-(samp:sample)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH DISTINCT ss, samp, collect(DISTINCT samp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
WHERE ss.disease_subtype IN ["Tubular Carcinoma"] and d.tumor_size_group In ["(3,4]"] and d.er_status In ["Positive"]and d.pr_status In ["Positive"] 
return DISTINCT ss.study_subject_id as `Case ID`,
   p.program_acronym as `Program Code`,
    p.program_id as Program_ID,
   s.study_acronym as `Arm`,
   ss.disease_subtype as `Diagnosis`,
   sf.grouped_recurrence_score AS `Recurrence Score`,
   d.tumor_size_group AS `tumor_size`,
   d.er_status AS `ER Status`,
   d.pr_status AS `PR Status`,
   demo.age_at_index AS `Age (years)`,
	demo.survival_time AS `Survival (days)`
'@
$ws.Range("B2").Value = $newQuery

# --- C2 / D2: updated deliverable file names ---
$ws.Range("C2").Value = "TC02_Bento_E2E_Select-Single-Add-To-Cart_Manifest.xlsx"
$ws.Range("D2").Value = "TC02_Bento_E2E_Select-Single-Add-To-Cart_WebData.xlsx"

# --- A2: updated QA url text + re-pointed hyperlink (drop stale #/ sub-address/display) ---
$ws.Range("A2").Value = "https://bento-qa.bento-tools.org/"
$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://caninecommons.cancer.gov/")
$ws.Range("A2").Style = "Hyperlink"

# --- Row 2 height nudged to the refreshed layout ---
$ws.Rows.Item(2).RowHeight = 375

# --- View: scroll so column A is back in view (topLeftCell A2) while keeping the B2 selection ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 2

Write-Host "edit applied"
